$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 2844.5
$ws.Range("I80").Value = 660.3333
$ws.Range("J80").Value = 4482.625
$ws.Range("K80").Value = 1980.9999
$ws.Range("L80").Value = 13447.875
$ws.Range("M80").Value = -982.9999
$ws.Range("N80").Value = -15443.875

$ws.Range("H83").Value = 2844.5
$ws.Range("I83").Value = 660.3333
$ws.Range("J83").Value = 4482.625
$ws.Range("K83").Value = 5942.9997
$ws.Range("L83").Value = 40343.625
$ws.Range("M83").Value = -950.9997000000003
$ws.Range("N83").Value = -50327.625

$ws.Range("H92").Value = 6250350
$ws.Range("I92").Value = 7143186
$ws.Range("K92").Value = 7143186
$ws.Range("M92").Value = -7141938

$ws.Range("H111").Value = 666
$ws.Range("J111").Value = 885
$ws.Range("L111").Value = 2655
$ws.Range("N111").Value = -8789

$ws.Range("H112").Value = 1318.625
$ws.Range("I112").Value = 1162.25
$ws.Range("J112").Value = 1349.9
$ws.Range("K112").Value = 3486.75
$ws.Range("L112").Value = 4049.7
$ws.Range("M112").Value = -2378.75
$ws.Range("N112").Value = -6265.700000000001

$ws.Range("H138").Value = 1793.6316
$ws.Range("I138").Value = 1710.5294
$ws.Range("J138").Value = 2500
$ws.Range("K138").Value = 5131.5882
$ws.Range("L138").Value = 7500
$ws.Range("M138").Value = 8.411799999999857
$ws.Range("N138").Value = -17780

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 590.1818
$ws.Range("I2").Value = 561.5
$ws.Range("J2").Value = 666.6667
$ws.Range("K2").Value = 561.5
$ws.Range("L2").Value = 666.6667
$ws.Range("M2").Value = -448.5
$ws.Range("N2").Value = -892.6667

$ws.Range("H32").Value = 7168.08
$ws.Range("I32").Value = 3419.3
$ws.Range("K32").Value = 3419.3
$ws.Range("M32").Value = -3132.3

$ws.Range("H102").Value = 46329
$ws.Range("I102").Value = 49641.19
$ws.Range("K102").Value = 49641.19
$ws.Range("M102").Value = -48019.19

$ws.Range("H116").Value = 590.1818
$ws.Range("I116").Value = 561.5
$ws.Range("J116").Value = 666.6667
$ws.Range("K116").Value = 561.5
$ws.Range("L116").Value = 666.6667
$ws.Range("M116").Value = 1732.5
$ws.Range("N116").Value = -5254.6667

$ws.Range("H121").Value = 48667.145
$ws.Range("J121").Value = 48667.145
$ws.Range("L121").Value = 48667.145
$ws.Range("N121").Value = -52161.145

$ws.Range("H132").Value = 2203
$ws.Range("I132").Value = 2143.7058
$ws.Range("J132").Value = 2539
$ws.Range("K132").Value = 6431.117400000001
$ws.Range("L132").Value = 7617
$ws.Range("M132").Value = -3901.117400000001
$ws.Range("N132").Value = -12677

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 590.1818
$ws.Range("I3").Value = 561.5
$ws.Range("J3").Value = 666.6667
$ws.Range("K3").Value = 561.5
$ws.Range("L3").Value = 666.6667
$ws.Range("M3").Value = -447.5
$ws.Range("N3").Value = -894.6667

$ws.Range("H20").Value = 1301.2307
$ws.Range("I20").Value = 1242.75
$ws.Range("K20").Value = 1242.75
$ws.Range("M20").Value = -995.75

$ws.Range("H132").Value = 29955.111
$ws.Range("J132").Value = 29955.111
$ws.Range("L132").Value = 29955.111
$ws.Range("N132").Value = -40075.111

$ws.Range("H134").Value = 6639.4116
$ws.Range("I134").Value = 4215.4546
$ws.Range("K134").Value = 12646.3638
$ws.Range("M134").Value = -10111.3638

$ws.Range("H138").Value = 99998.336
$ws.Range("J138").Value = 99998.336
$ws.Range("L138").Value = 99998.336
$ws.Range("N138").Value = -110278.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2398.3333
$ws.Range("I31").Value = 1583.3182
$ws.Range("K31").Value = 1583.3182
$ws.Range("M31").Value = -1288.3182

$ws.Range("H34").Value = 2398.3333
$ws.Range("I34").Value = 1583.3182
$ws.Range("K34").Value = 1583.3182
$ws.Range("M34").Value = -1381.3182

$ws.Range("H57").Value = 17250
$ws.Range("J57").Value = 21333.334
$ws.Range("L57").Value = 21333.334
$ws.Range("N57").Value = -22453.334

$ws.Range("H107").Value = 1402.2727
$ws.Range("I107").Value = 880.55554
$ws.Range("J107").Value = 3750
$ws.Range("K107").Value = 880.55554
$ws.Range("L107").Value = 3750
$ws.Range("M107").Value = 1039.44446
$ws.Range("N107").Value = -7590

$ws.Range("H134").Value = 1726651.9
$ws.Range("I134").Value = 2103026.8
$ws.Range("J134").Value = 127058.875
$ws.Range("K134").Value = 6309080.399999999
$ws.Range("L134").Value = 381176.625
$ws.Range("M134").Value = -6306545.399999999
$ws.Range("N134").Value = -386246.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 5473.273
$ws.Range("I56").Value = 5473.273
$ws.Range("K56").Value = 5473.273
$ws.Range("M56").Value = -4943.273

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 20000
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 20000
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H102").Value = 1978.2
$ws.Range("I102").Value = 1635.375
$ws.Range("K102").Value = 1635.375
$ws.Range("M102").Value = -13.375

$ws.Range("H122").Value = 4158.1
$ws.Range("J122").Value = 3583.1667
$ws.Range("L122").Value = 10749.5001
$ws.Range("N122").Value = -15649.5001

$ws.Range("H132").Value = 5201.3
$ws.Range("I132").Value = 1949.5
$ws.Range("J132").Value = 6014.25
$ws.Range("K132").Value = 5848.5
$ws.Range("L132").Value = 18042.75
$ws.Range("M132").Value = -3318.5
$ws.Range("N132").Value = -23102.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3942.9
$ws.Range("I16").Value = 3341.125
$ws.Range("J16").Value = 6350
$ws.Range("K16").Value = 3341.125
$ws.Range("L16").Value = 6350
$ws.Range("M16").Value = -3171.125
$ws.Range("N16").Value = -6690

$ws.Range("H22").Value = 2340.818
$ws.Range("I22").Value = 1791.6666
$ws.Range("K22").Value = 1791.6666
$ws.Range("M22").Value = -1496.6666

$ws.Range("H27").Value = 2340.818
$ws.Range("I27").Value = 1791.6666
$ws.Range("K27").Value = 1791.6666
$ws.Range("M27").Value = -1684.6666

$ws.Range("H46").Value = 2426.7144
$ws.Range("I46").Value = 1734.7142
$ws.Range("J46").Value = 3118.7144
$ws.Range("K46").Value = 1734.7142
$ws.Range("L46").Value = 3118.7144
$ws.Range("M46").Value = -1546.7142
$ws.Range("N46").Value = -3494.7144

$ws.Range("H47").Value = 32500
$ws.Range("J47").Value = 32500
$ws.Range("L47").Value = 32500
$ws.Range("N47").Value = -33480

$ws.Range("H52").Value = 32500
$ws.Range("J52").Value = 32500
$ws.Range("L52").Value = 32500
$ws.Range("N52").Value = -32966

$ws.Range("H61").Value = 1275
$ws.Range("I61").Value = 1275
$ws.Range("K61").Value = 1275
$ws.Range("M61").Value = -1073

$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").ClearContents()
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = 0

$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").ClearContents()
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = 0

$ws.Range("H68").Value = 5750.5
$ws.Range("I68").Value = 5750.5
$ws.Range("K68").Value = 5750.5
$ws.Range("M68").Value = -5001.5

$ws.Range("H71").Value = 5750.5
$ws.Range("I71").Value = 5750.5
$ws.Range("K71").Value = 28752.5
$ws.Range("M71").Value = -25008.5

$ws.Range("H82").Value = 2583
$ws.Range("I82").Value = 2041.3334
$ws.Range("K82").Value = 2041.3334
$ws.Range("M82").Value = -1680.3334

$ws.Range("H85").Value = 2583
$ws.Range("I85").Value = 2041.3334
$ws.Range("K85").Value = 2041.3334
$ws.Range("M85").Value = -793.3334

$ws.Range("H113").Value = 1275
$ws.Range("I113").Value = 1275
$ws.Range("K113").Value = 1275
$ws.Range("M113").Value = 895

$ws.Range("H119").Value = 29710
$ws.Range("J119").Value = 29710
$ws.Range("L119").Value = 29710
$ws.Range("N119").Value = -39386

$ws.Range("H122").Value = 171432850
$ws.Range("I122").Value = 200004670
$ws.Range("J122").Value = 100003250
$ws.Range("K122").Value = 600014010
$ws.Range("L122").Value = 300009750
$ws.Range("M122").Value = -600011560
$ws.Range("N122").Value = -300014650

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1677.6578
$ws.Range("I132").Value = 1244.5385
$ws.Range("J132").Value = 2616.0833
$ws.Range("K132").Value = 3733.6155
$ws.Range("L132").Value = 7848.249899999999
$ws.Range("M132").Value = -1203.6155
$ws.Range("N132").Value = -12908.2499

$ws.Range("H136").Value = 1457.3334
$ws.Range("I136").Value = 1243.4
$ws.Range("J136").Value = 1724.75
$ws.Range("K136").Value = 3730.2
$ws.Range("L136").Value = 5174.25
$ws.Range("M136").Value = -1180.2
$ws.Range("N136").Value = -10274.25
